$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) to the data-revenue table, copying the
# formatting from the existing last year column (R) for the header (row 4)
# and the data value (row 5).
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 76.099999999999994

$ws.Range("P8").Select()
